# B6-PowerPoint.pptx edit:
#  1. Re-style the three summary tables (slides 14-16) from the deck's
#     default custom table style to the built-in "Medium Style 2 - Accent 1"
#     table style ({A156A0F2-606A-44A7-9B54-8AA12F1F39EB}).
#  2. Swap the theme colours so the deck's main theme becomes the stock
#     "Office" colour scheme (previously the "Integral"/"Red Violet"
#     scheme that shipped with the template).

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newStyleId = "{A156A0F2-606A-44A7-9B54-8AA12F1F39EB}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Theme colours --------------------------------------------------
# New ("Office") colour scheme values, in PpColorSchemeIndex order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
$officeColors = @{
    1  = 0x000000
    2  = 0xFFFFFF
    3  = 0x6A5444   # 44546A stored BGR
    4  = 0xE6E6E7   # E7E6E6 stored BGR
    5  = 0xD59B5B   # 5B9BD5 stored BGR
    6  = 0x317DED   # ED7D31 stored BGR
    7  = 0xA5A5A5   # A5A5A5 stored BGR
    8  = 0x00C0FF   # FFC000 stored BGR
    9  = 0xC47244   # 4472C4 stored BGR
    10 = 0x47AD70   # 70AD47 stored BGR
    11 = 0xC16305   # 0563C1 stored BGR
    12 = 0x724F95   # 954F72 stored BGR
}

$master = $p.SlideMaster
$scheme = $master.ColorScheme
foreach ($idx in 1..12) {
    $scheme.Colors($idx).RGB = $officeColors[$idx]
}
